$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CorrCell($row, $col, $expected, $value) {
    $cell = $t.Cell($row, $col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $expected) {
        Write-Output "WARNING: cell ($row,$col) expected [$expected] but found [$current]"
    }
    $cell.Range.Text = $value
}

Set-CorrCell 2 2 "0.26" "0.29"
Set-CorrCell 2 5 "-0.04" "-0.03"
Set-CorrCell 2 6 "0.05" "0.04"
Set-CorrCell 3 1 "0.26" "0.29"
Set-CorrCell 3 3 "0.14" "0.15"
Set-CorrCell 3 4 "-0.22" "-0.21"
Set-CorrCell 3 5 "-0.22" "-0.24"
Set-CorrCell 3 6 "0.13" "0.12"
Set-CorrCell 3 7 "-0.14" "-0.15"
Set-CorrCell 4 2 "0.14" "0.15"
Set-CorrCell 4 4 "-0.08" "-0.07"
Set-CorrCell 4 5 "0.07" "0.1"
Set-CorrCell 4 6 "-0.04" "-0.05"
Set-CorrCell 5 2 "-0.22" "-0.21"
Set-CorrCell 5 3 "-0.08" "-0.07"
Set-CorrCell 5 5 "0.03" "0.02"
Set-CorrCell 5 7 "0.08" "0.07"
Set-CorrCell 6 1 "-0.04" "-0.03"
Set-CorrCell 6 2 "-0.22" "-0.24"
Set-CorrCell 6 3 "0.07" "0.1"
Set-CorrCell 6 4 "0.03" "0.02"
Set-CorrCell 6 7 "-0.07" "-0.1"
Set-CorrCell 7 1 "0.05" "0.04"
Set-CorrCell 7 2 "0.13" "0.12"
Set-CorrCell 7 3 "-0.04" "-0.05"
Set-CorrCell 7 7 "0.04" "0.05"
Set-CorrCell 8 2 "-0.14" "-0.15"
Set-CorrCell 8 4 "0.08" "0.07"
Set-CorrCell 8 5 "-0.07" "-0.1"
Set-CorrCell 8 6 "0.04" "0.05"

Write-Output "done"
